$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 230
$ws.Range("A230").Value() = 7
$ws.Range("B230").Value() = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C230").Value() = "Ñuble"
$ws.Range("D230").Value() = 44595
$ws.Range("D230").NumberFormat() = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E230").Value() = 16
$ws.Range("F230").Value() = "Fruta"
$ws.Range("G230").Value() = 100103
$ws.Range("H230").Value() = "Frutos de hueso (carozo)"
$ws.Range("I230").Value() = 100103004
$ws.Range("J230").Value() = "Durazno"
$ws.Range("K230").Value() = "Carson"
$ws.Range("L230").Value() = "Especial"
$ws.Range("M230").Value() = 80
$ws.Range("N230").Value() = 12000
$ws.Range("O230").Value() = 12000
$ws.Range("P230").Value() = 12000
$ws.Range("Q230").Value() = "`$/caja 16 kilos empedrada"
$ws.Range("R230").Value() = "Región de O'Higgins"
$ws.Range("S230").Value() = 750
$ws.Range("T230").Value() = 16

# Row 231
$ws.Range("A231").Value() = 7
$ws.Range("B231").Value() = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C231").Value() = "Ñuble"
$ws.Range("D231").Value() = 44595
$ws.Range("D231").NumberFormat() = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E231").Value() = 16
$ws.Range("F231").Value() = "Fruta"
$ws.Range("G231").Value() = 100103
$ws.Range("H231").Value() = "Frutos de hueso (carozo)"
$ws.Range("I231").Value() = 100103004
$ws.Range("J231").Value() = "Durazno"
$ws.Range("K231").Value() = "Carson"
$ws.Range("L231").Value() = "Primera"
$ws.Range("M231").Value() = 120
$ws.Range("N231").Value() = 10000
$ws.Range("O231").Value() = 11000
$ws.Range("P231").Value() = 10500
$ws.Range("Q231").Value() = "`$/caja 16 kilos empedrada"
$ws.Range("R231").Value() = "Región de O'Higgins"
$ws.Range("S231").Value() = 656
$ws.Range("T231").Value() = 16

# Row 232
$ws.Range("A232").Value() = 7
$ws.Range("B232").Value() = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C232").Value() = "Ñuble"
$ws.Range("D232").Value() = 44595
$ws.Range("D232").NumberFormat() = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E232").Value() = 16
$ws.Range("F232").Value() = "Fruta"
$ws.Range("G232").Value() = 100103
$ws.Range("H232").Value() = "Frutos de hueso (carozo)"
$ws.Range("I232").Value() = 100103004
$ws.Range("J232").Value() = "Durazno"
$ws.Range("K232").Value() = "Elegant Lady"
$ws.Range("L232").Value() = "Especial"
$ws.Range("M232").Value() = 80
$ws.Range("N232").Value() = 12000
$ws.Range("O232").Value() = 12000
$ws.Range("P232").Value() = 12000
$ws.Range("Q232").Value() = "`$/caja 16 kilos empedrada"
$ws.Range("R232").Value() = "Región de O'Higgins"
$ws.Range("S232").Value() = 750
$ws.Range("T232").Value() = 16

# Row 233
$ws.Range("A233").Value() = 7
$ws.Range("B233").Value() = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C233").Value() = "Ñuble"
$ws.Range("D233").Value() = 44595
$ws.Range("D233").NumberFormat() = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E233").Value() = 16
$ws.Range("F233").Value() = "Fruta"
$ws.Range("G233").Value() = 100103
$ws.Range("H233").Value() = "Frutos de hueso (carozo)"
$ws.Range("I233").Value() = 100103004
$ws.Range("J233").Value() = "Durazno"
$ws.Range("K233").Value() = "Elegant Lady"
$ws.Range("L233").Value() = "Primera"
$ws.Range("M233").Value() = 120
$ws.Range("N233").Value() = 10000
$ws.Range("O233").Value() = 11000
$ws.Range("P233").Value() = 10500
$ws.Range("Q233").Value() = "`$/caja 16 kilos empedrada"
$ws.Range("R233").Value() = "Región de O'Higgins"
$ws.Range("S233").Value() = 656
$ws.Range("T233").Value() = 16
